$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("G4").Value2 = $ws.Range("G4").Value2 -replace [regex]::Escape("link_Click(VT200_0576_mainpage_xpath)"), "ClickNativeIcon(VT200_0576_mainpage_xpath)"
$ws.Range("G6").Value2 = $ws.Range("G6").Value2 -replace [regex]::Escape("link_Click(VT200_0576_mainpage_xpath)"), "ClickNativeIcon(VT200_0576_mainpage_xpath)"
$ws.Range("G7").Value2 = $ws.Range("G7").Value2 -replace [regex]::Escape("link_Click(VT200_0576_mainpage_xpath)"), "ClickNativeIcon(VT200_0576_mainpage_xpath)"
$ws.Range("G11").Value2 = $ws.Range("G11").Value2 -replace [regex]::Escape("link_Click(VT200_0585_page1_xpath)"), "ClickNativeIcon(VT200_0585_page1_xpath)"

$ws.Range("J2:J13").Select()
